$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 73.38544233333333
$ws.Range("H2").Value = 220.156327
$ws.Range("I2").Value = 0.1214979676060253
$ws.Range("J2").Value = 0.1214979676060253
$ws.Range("M2").Value = 16.18329966666667
$ws.Range("N2").Value = 48.549899
$ws.Range("O2").Value = 0.6688494993193013
$ws.Range("P2").Value = 0.6688494993193014
$ws.Range("Q2").Value = 1187.618604451219
$ws.Range("R2").Value = 10688.56744006097
$ws.Range("S2").Value = 0.08126385480160271
$ws.Range("T2").Value = 0.08126385480160271

$ws.Range("G3").Value = 73.38544233333333
$ws.Range("H3").Value = 220.156327
$ws.Range("I3").Value = 0.1214979676060253
$ws.Range("J3").Value = 0.1214979676060253
$ws.Range("M3").Value = 0.2227456666666666
$ws.Range("N3").Value = 0.668237
$ws.Range("O3").Value = 0.00920599202228272
$ws.Range("P3").Value = 0.009205992022282724
$ws.Range("Q3").Value = 16.34628927616656
$ws.Range("R3").Value = 147.116603485499
$ws.Range("S3").Value = 0.001118509320504633
$ws.Range("T3").Value = 0.001118509320504634

$ws.Range("G4").Value = 73.38544233333333
$ws.Range("H4").Value = 220.156327
$ws.Range("I4").Value = 0.1214979676060253
$ws.Range("J4").Value = 0.1214979676060253
$ws.Range("M4").Value = 3.563453666666667
$ws.Range("N4").Value = 10.690361
$ws.Range("O4").Value = 0.1472761581315048
$ws.Range("P4").Value = 0.1472761581315048
$ws.Range("Q4").Value = 261.5056235626719
$ws.Range("R4").Value = 2353.550612064047
$ws.Range("S4").Value = 0.01789375388980143
$ws.Range("T4").Value = 0.01789375388980143

$ws.Range("G5").Value = 73.38544233333333
$ws.Range("H5").Value = 220.156327
$ws.Range("I5").Value = 0.1214979676060253
$ws.Range("J5").Value = 0.1214979676060253
$ws.Range("M5").Value = 0.4040503333333333
$ws.Range("N5").Value = 1.212151
$ws.Range("O5").Value = 0.01669924358543754
$ws.Range("P5").Value = 0.01669924358543754
$ws.Range("Q5").Value = 29.65141243659744
$ws.Range("R5").Value = 266.862711929377
$ws.Range("S5").Value = 0.002028924156188616
$ws.Range("T5").Value = 0.002028924156188616

$ws.Range("G6").Value = 73.38544233333333
$ws.Range("H6").Value = 220.156327
$ws.Range("I6").Value = 0.1214979676060253
$ws.Range("J6").Value = 0.1214979676060253
$ws.Range("M6").Value = 3.822177333333334
$ws.Range("N6").Value = 11.466532
$ws.Range("O6").Value = 0.1579691069414737
$ws.Range("P6").Value = 0.1579691069414737
$ws.Range("Q6").Value = 280.4921742831071
$ws.Range("R6").Value = 2524.429568547964
$ws.Range("S6").Value = 0.01919292543792792
$ws.Range("T6").Value = 0.01919292543792792

$ws.Range("I7").Value = 0.3924995450689984
$ws.Range("J7").Value = 0.3924995450689983
$ws.Range("M7").Value = 16.18329966666667
$ws.Range("N7").Value = 48.549899
$ws.Range("O7").Value = 0.6688494993193013
$ws.Range("P7").Value = 0.6688494993193014
$ws.Range("Q7").Value = 3836.605427624171
$ws.Range("R7").Value = 34529.44884861753
$ws.Range("S7").Value = 0.2625231242024531
$ws.Range("T7").Value = 0.2625231242024531

$ws.Range("I8").Value = 0.3924995450689984
$ws.Range("J8").Value = 0.3924995450689983
$ws.Range("M8").Value = 0.2227456666666666
$ws.Range("N8").Value = 0.668237
$ws.Range("O8").Value = 0.00920599202228272
$ws.Range("P8").Value = 0.009205992022282724
$ws.Range("Q8").Value = 52.80673603747956
$ws.Range("R8").Value = 475.260624337316
$ws.Range("S8").Value = 0.003613347680654797
$ws.Range("T8").Value = 0.003613347680654797

$ws.Range("I9").Value = 0.3924995450689984
$ws.Range("J9").Value = 0.3924995450689983
$ws.Range("M9").Value = 3.563453666666667
$ws.Range("N9").Value = 10.690361
$ws.Range("O9").Value = 0.1472761581315048
$ws.Range("P9").Value = 0.1472761581315048
$ws.Range("Q9").Value = 844.7946933084611
$ws.Range("R9").Value = 7603.152239776149
$ws.Range("S9").Value = 0.05780582506612549
$ws.Range("T9").Value = 0.05780582506612549

$ws.Range("I10").Value = 0.3924995450689984
$ws.Range("J10").Value = 0.3924995450689983
$ws.Range("M10").Value = 0.4040503333333333
$ws.Range("N10").Value = 1.212151
$ws.Range("O10").Value = 0.01669924358543754
$ws.Range("P10").Value = 0.01669924358543754
$ws.Range("Q10").Value = 95.78897590909646
$ws.Range("R10").Value = 862.100783181868
$ws.Range("S10").Value = 0.006554445510280624
$ws.Range("T10").Value = 0.006554445510280623

$ws.Range("I11").Value = 0.3924995450689984
$ws.Range("J11").Value = 0.3924995450689983
$ws.Range("M11").Value = 3.822177333333334
$ws.Range("N11").Value = 11.466532
$ws.Range("O11").Value = 0.1579691069414737
$ws.Range("P11").Value = 0.1579691069414737
$ws.Range("Q11").Value = 906.1308017803754
$ws.Range("R11").Value = 8155.177216023377
$ws.Range("S11").Value = 0.06200280260948439
$ws.Range("T11").Value = 0.06200280260948438

$ws.Range("G12").Value = 138.1628113333333
$ws.Range("H12").Value = 414.488434
$ws.Range("I12").Value = 0.2287442882675098
$ws.Range("J12").Value = 0.2287442882675098
$ws.Range("M12").Value = 16.18329966666667
$ws.Range("N12").Value = 48.549899
$ws.Range("O12").Value = 0.6688494993193013
$ws.Range("P12").Value = 0.6688494993193014
$ws.Range("Q12").Value = 2235.930178596463
$ws.Range("R12").Value = 20123.37160736816
$ws.Range("S12").Value = 0.1529955026798739
$ws.Range("T12").Value = 0.1529955026798739

$ws.Range("G13").Value = 138.1628113333333
$ws.Range("H13").Value = 414.488434
$ws.Range("I13").Value = 0.2287442882675098
$ws.Range("J13").Value = 0.2287442882675098
$ws.Range("M13").Value = 0.2227456666666666
$ws.Range("N13").Value = 0.668237
$ws.Range("O13").Value = 0.00920599202228272
$ws.Range("P13").Value = 0.009205992022282724
$ws.Range("Q13").Value = 30.77516751898422
$ws.Range("R13").Value = 276.976507670858
$ws.Range("S13").Value = 0.002105818092933435
$ws.Range("T13").Value = 0.002105818092933435

$ws.Range("G14").Value = 138.1628113333333
$ws.Range("H14").Value = 414.488434
$ws.Range("I14").Value = 0.2287442882675098
$ws.Range("J14").Value = 0.2287442882675098
$ws.Range("M14").Value = 3.563453666666667
$ws.Range("N14").Value = 10.690361
$ws.Range("O14").Value = 0.1472761581315048
$ws.Range("P14").Value = 0.1472761581315048
$ws.Range("Q14").Value = 492.3367766427417
$ws.Range("R14").Value = 4431.030989784675
$ws.Range("S14").Value = 0.03368857997056429
$ws.Range("T14").Value = 0.03368857997056428

$ws.Range("G15").Value = 138.1628113333333
$ws.Range("H15").Value = 414.488434
$ws.Range("I15").Value = 0.2287442882675098
$ws.Range("J15").Value = 0.2287442882675098
$ws.Range("M15").Value = 0.4040503333333333
$ws.Range("N15").Value = 1.212151
$ws.Range("O15").Value = 0.01669924358543754
$ws.Range("P15").Value = 0.01669924358543754
$ws.Range("Q15").Value = 55.82472997350378
$ws.Range("R15").Value = 502.422569761534
$ws.Range("S15").Value = 0.003819856588556689
$ws.Range("T15").Value = 0.003819856588556689

$ws.Range("G16").Value = 138.1628113333333
$ws.Range("H16").Value = 414.488434
$ws.Range("I16").Value = 0.2287442882675098
$ws.Range("J16").Value = 0.2287442882675098
$ws.Range("M16").Value = 3.822177333333334
$ws.Range("N16").Value = 11.466532
$ws.Range("O16").Value = 0.1579691069414737
$ws.Range("P16").Value = 0.1579691069414737
$ws.Range("Q16").Value = 528.0827657878765
$ws.Range("R16").Value = 4752.744892090888
$ws.Range("S16").Value = 0.03613453093558155
$ws.Range("T16").Value = 0.03613453093558155

$ws.Range("G17").Value = 49.051656
$ws.Range("H17").Value = 147.154968
$ws.Range("I17").Value = 0.08121060965524597
$ws.Range("J17").Value = 0.08121060965524596
$ws.Range("M17").Value = 16.18329966666667
$ws.Range("N17").Value = 48.549899
$ws.Range("O17").Value = 0.6688494993193013
$ws.Range("P17").Value = 0.6688494993193014
$ws.Range("Q17").Value = 793.817648194248
$ws.Range("R17").Value = 7144.358833748232
$ws.Range("S17").Value = 0.05431767560732648
$ws.Range("T17").Value = 0.05431767560732648

$ws.Range("G18").Value = 49.051656
$ws.Range("H18").Value = 147.154968
$ws.Range("I18").Value = 0.08121060965524597
$ws.Range("J18").Value = 0.08121060965524596
$ws.Range("M18").Value = 0.2227456666666666
$ws.Range("N18").Value = 0.668237
$ws.Range("O18").Value = 0.00920599202228272
$ws.Range("P18").Value = 0.009205992022282724
$ws.Range("Q18").Value = 10.926043816824
$ws.Range("R18").Value = 98.334394351416
$ws.Range("S18").Value = 0.0007476242246109105
$ws.Range("T18").Value = 0.0007476242246109106

$ws.Range("G19").Value = 49.051656
$ws.Range("H19").Value = 147.154968
$ws.Range("I19").Value = 0.08121060965524597
$ws.Range("J19").Value = 0.08121060965524596
$ws.Range("M19").Value = 3.563453666666667
$ws.Range("N19").Value = 10.690361
$ws.Range("O19").Value = 0.1472761581315048
$ws.Range("P19").Value = 0.1472761581315048
$ws.Range("Q19").Value = 174.793303429272
$ws.Range("R19").Value = 1573.139730863448
$ws.Range("S19").Value = 0.01196038658954191
$ws.Range("T19").Value = 0.01196038658954191

$ws.Range("G20").Value = 49.051656
$ws.Range("H20").Value = 147.154968
$ws.Range("I20").Value = 0.08121060965524597
$ws.Range("J20").Value = 0.08121060965524596
$ws.Range("M20").Value = 0.4040503333333333
$ws.Range("N20").Value = 1.212151
$ws.Range("O20").Value = 0.01669924358543754
$ws.Range("P20").Value = 0.01669924358543754
$ws.Range("Q20").Value = 19.819337957352
$ws.Range("R20").Value = 178.374041616168
$ws.Range("S20").Value = 0.001356155752354838
$ws.Range("T20").Value = 0.001356155752354838

$ws.Range("G21").Value = 49.051656
$ws.Range("H21").Value = 147.154968
$ws.Range("I21").Value = 0.08121060965524597
$ws.Range("J21").Value = 0.08121060965524596
$ws.Range("M21").Value = 3.822177333333334
$ws.Range("N21").Value = 11.466532
$ws.Range("O21").Value = 0.1579691069414737
$ws.Range("P21").Value = 0.1579691069414737
$ws.Range("Q21").Value = 187.484127725664
$ws.Range("R21").Value = 1687.357149530976
$ws.Range("S21").Value = 0.01282876748141183
$ws.Range("T21").Value = 0.01282876748141183

$ws.Range("G22").Value = 106.3337146666667
$ws.Range("H22").Value = 319.001144
$ws.Range("I22").Value = 0.1760475894022206
$ws.Range("J22").Value = 0.1760475894022206
$ws.Range("M22").Value = 16.18329966666667
$ws.Range("N22").Value = 48.549899
$ws.Range("O22").Value = 0.6688494993193013
$ws.Range("P22").Value = 0.6688494993193014
$ws.Range("Q22").Value = 1720.830369120495
$ws.Range("R22").Value = 15487.47332208446
$ws.Range("S22").Value = 0.1177493420280452
$ws.Range("T22").Value = 0.1177493420280452

$ws.Range("G23").Value = 106.3337146666667
$ws.Range("H23").Value = 319.001144
$ws.Range("I23").Value = 0.1760475894022206
$ws.Range("J23").Value = 0.1760475894022206
$ws.Range("M23").Value = 0.2227456666666666
$ws.Range("N23").Value = 0.668237
$ws.Range("O23").Value = 0.00920599202228272
$ws.Range("P23").Value = 0.009205992022282724
$ws.Range("Q23").Value = 23.68537416256978
$ws.Range("R23").Value = 213.168367463128
$ws.Range("S23").Value = 0.001620692703578947
$ws.Range("T23").Value = 0.001620692703578947

$ws.Range("G24").Value = 106.3337146666667
$ws.Range("H24").Value = 319.001144
$ws.Range("I24").Value = 0.1760475894022206
$ws.Range("J24").Value = 0.1760475894022206
$ws.Range("M24").Value = 3.563453666666667
$ws.Range("N24").Value = 10.690361
$ws.Range("O24").Value = 0.1472761581315048
$ws.Range("P24").Value = 0.1472761581315048
$ws.Range("Q24").Value = 378.9152654192205
$ws.Range("R24").Value = 3410.237388772985
$ws.Range("S24").Value = 0.02592761261547166
$ws.Range("T24").Value = 0.02592761261547167

$ws.Range("G25").Value = 106.3337146666667
$ws.Range("H25").Value = 319.001144
$ws.Range("I25").Value = 0.1760475894022206
$ws.Range("J25").Value = 0.1760475894022206
$ws.Range("M25").Value = 0.4040503333333333
$ws.Range("N25").Value = 1.212151
$ws.Range("O25").Value = 0.01669924358543754
$ws.Range("P25").Value = 0.01669924358543754
$ws.Range("Q25").Value = 42.96417285563822
$ws.Range("R25").Value = 386.677555700744
$ws.Range("S25").Value = 0.002939861578056774
$ws.Range("T25").Value = 0.002939861578056774

$ws.Range("G26").Value = 106.3337146666667
$ws.Range("H26").Value = 319.001144
$ws.Range("I26").Value = 0.1760475894022206
$ws.Range("J26").Value = 0.1760475894022206
$ws.Range("M26").Value = 3.822177333333334
$ws.Range("N26").Value = 11.466532
$ws.Range("O26").Value = 0.1579691069414737
$ws.Range("P26").Value = 0.1579691069414737
$ws.Range("Q26").Value = 406.4263139680676
$ws.Range("R26").Value = 3657.836825712609
$ws.Range("S26").Value = 0.02781008047706804
$ws.Range("T26").Value = 0.02781008047706804
